# Update "想去人数" (interest count) values in column F across the four
# sheets to reflect the newer snapshot of generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1272
$ws1.Range("F5").Value  = 337
$ws1.Range("F6").Value  = 3852
$ws1.Range("F8").Value  = 754
$ws1.Range("F9").Value  = 2241
$ws1.Range("F11").Value = 221
$ws1.Range("F12").Value = 740
$ws1.Range("F13").Value = 166
$ws1.Range("F14").Value = 166
$ws1.Range("F15").Value = 2161
$ws1.Range("F17").Value = 10
$ws1.Range("F18").Value = 49
$ws1.Range("F19").Value = 338
$ws1.Range("F20").Value = 227
$ws1.Range("F21").Value = 28
$ws1.Range("F22").Value = 270

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 44
$ws2.Range("F7").Value  = 11
$ws2.Range("F8").Value  = 128
$ws2.Range("F9").Value  = 95
$ws2.Range("F11").Value = 88
$ws2.Range("F12").Value = 226
$ws2.Range("F16").Value = 8

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6397
$ws3.Range("F3").Value = 822
$ws3.Range("F4").Value = 2087
$ws3.Range("F5").Value = 318

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6397
$ws4.Range("F3").Value  = 822
$ws4.Range("F4").Value  = 2087
$ws4.Range("F5").Value  = 318
$ws4.Range("F7").Value  = 44
$ws4.Range("F8").Value  = 44
$ws4.Range("F11").Value = 1272
$ws4.Range("F15").Value = 11
$ws4.Range("F16").Value = 337
$ws4.Range("F17").Value = 3852
$ws4.Range("F18").Value = 128
$ws4.Range("F20").Value = 95
$ws4.Range("F22").Value = 88
$ws4.Range("F23").Value = 754
$ws4.Range("F24").Value = 2241
$ws4.Range("F26").Value = 226
$ws4.Range("F27").Value = 221
$ws4.Range("F28").Value = 740
$ws4.Range("F29").Value = 166
$ws4.Range("F30").Value = 166
$ws4.Range("F32").Value = 2161
$ws4.Range("F36").Value = 10
$ws4.Range("F37").Value = 49
$ws4.Range("F38").Value = 338
$ws4.Range("F39").Value = 227
$ws4.Range("F40").Value = 28
$ws4.Range("F41").Value = 8
$ws4.Range("F48").Value = 270
